# Implement new Excel-based structure with en/zh folders and remove
# translation functionality.
#
# The "en-nz" / "zh-nz" webpage rows (rows 12-34, currently in column A)
# move over to column H, leaving the en-au/zh-au rows (7-11 and the
# header block 1-6) untouched in column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 12; $r -le 34; $r++) {
    $ws.Range("A$r").Cut($ws.Range("H$r"))
    $ws.Range("A$r").Clear()
}

# Match the saved selection: H12:H34 selected, active cell H12.
[void]$ws.Range("H12:H34").Select()
